$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 5657
$ws.Range("I3").Value = 5918
$ws.Range("I5").Value = 540
$ws.Range("I6").Value = 6611
$ws.Range("I7").Value = 20066

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I2").Value = 158
$ws.Range("I4").Value = 79
$ws.Range("I5").Value = 65
$ws.Range("I7").Value = 634
$ws.Range("I8").Value = 1203
$ws.Range("I11").Value = 296
$ws.Range("I18").Value = 143
$ws.Range("I19").Value = 551
$ws.Range("I20").Value = 482
$ws.Range("I21").Value = 91
$ws.Range("I29").Value = 1257
$ws.Range("I33").Value = 913
$ws.Range("I34").Value = 96
$ws.Range("I36").Value = 261
$ws.Range("I37").Value = 642
$ws.Range("I42").Value = 673
$ws.Range("I48").Value = 270
$ws.Range("I49").Value = 134
$ws.Range("I50").Value = 97
$ws.Range("I52").Value = 439
$ws.Range("I53").Value = 210
$ws.Range("I54").Value = 421
$ws.Range("I55").Value = 223
$ws.Range("I60").Value = 110
$ws.Range("I63").Value = 71
$ws.Range("I65").Value = 469
$ws.Range("I67").Value = 790
$ws.Range("I68").Value = 71
$ws.Range("I75").Value = 63
$ws.Range("I76").Value = 293
$ws.Range("I77").Value = 128
$ws.Range("I79").Value = 569
$ws.Range("I81").Value = 18
$ws.Range("I83").Value = 426
$ws.Range("I84").Value = 174
$ws.Range("I85").Value = 905
$ws.Range("I90").Value = 248
$ws.Range("I94").Value = 210
$ws.Range("I95").Value = 316
$ws.Range("I96").Value = 215
$ws.Range("I97").Value = 165
$ws.Range("I99").Value = 370
$ws.Range("I101").Value = 20066

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I3").Value = 351
$ws.Range("I6").Value = 232
$ws.Range("I7").Value = 905

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I3").Value = 158
$ws.Range("I7").Value = 439

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I6").Value = 80
$ws.Range("I7").Value = 296

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 372
$ws.Range("I3").Value = 337
$ws.Range("I6").Value = 390
$ws.Range("I7").Value = 1203

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("I2").Value = 46
$ws.Range("I3").Value = 45
$ws.Range("I7").Value = 210

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I6").Value = 166
$ws.Range("I7").Value = 634

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I2").Value = 64
$ws.Range("I5").Value = 3
$ws.Range("I7").Value = 215

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I2").Value = 197
$ws.Range("I3").Value = 213
$ws.Range("I6").Value = 183
$ws.Range("I7").Value = 642

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I2").Value = 102
$ws.Range("I3").Value = 136
$ws.Range("I7").Value = 370

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I3").Value = 289
$ws.Range("I6").Value = 248
$ws.Range("I7").Value = 790

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("I2").Value = 64
$ws.Range("I3").Value = 57
$ws.Range("I7").Value = 174

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I2").Value = 157
$ws.Range("I3").Value = 141
$ws.Range("I6").Value = 133
$ws.Range("I7").Value = 469

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I3").Value = 158
$ws.Range("I7").Value = 426

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I3").Value = 114
$ws.Range("I6").Value = 65
$ws.Range("I7").Value = 316

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 208
$ws.Range("I6").Value = 284
$ws.Range("I7").Value = 913

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("I6").Value = 82
$ws.Range("I7").Value = 134

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I4").Value = 29
$ws.Range("I6").Value = 199
$ws.Range("I7").Value = 421

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 368
$ws.Range("I3").Value = 435
$ws.Range("I6").Value = 344
$ws.Range("I7").Value = 1257

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I3").Value = 169
$ws.Range("I6").Value = 160
$ws.Range("I7").Value = 551

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("I4").Value = 31
$ws.Range("I7").Value = 270

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I6").Value = 136
$ws.Range("I7").Value = 293

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I3").Value = 230
$ws.Range("I6").Value = 201
$ws.Range("I7").Value = 673

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("I3").Value = 71
$ws.Range("I7").Value = 223

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("I6").Value = 67
$ws.Range("I7").Value = 91

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I2").Value = 164
$ws.Range("I3").Value = 181
$ws.Range("I6").Value = 166
$ws.Range("I7").Value = 569

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I2").Value = 138
$ws.Range("I6").Value = 160
$ws.Range("I7").Value = 482

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("I3").Value = 34
$ws.Range("I7").Value = 143

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I2").Value = 79
$ws.Range("I7").Value = 261

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("I6").Value = 20
$ws.Range("I7").Value = 96

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("I3").Value = 35
$ws.Range("I7").Value = 210

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("I2").Value = 28
$ws.Range("I7").Value = 97

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("I2").Value = 52
$ws.Range("I7").Value = 158

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("I2").Value = 27
$ws.Range("I6").Value = 103
$ws.Range("I7").Value = 165

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("I6").Value = 28
$ws.Range("I7").Value = 65

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("I6").Value = 15
$ws.Range("I7").Value = 63

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I3").Value = 62
$ws.Range("I6").Value = 82
$ws.Range("I7").Value = 248

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("I3").Value = 22
$ws.Range("I7").Value = 71

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("I6").Value = 33
$ws.Range("I7").Value = 110

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("I2").Value = 41
$ws.Range("I3").Value = 44
$ws.Range("I7").Value = 128

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("I2").Value = 31
$ws.Range("I7").Value = 79

$ws = $wb.Worksheets.Item('Sauganash,Forest Glen')
$ws.Range("I3").Value = 5
$ws.Range("I6").Value = 18
